$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (bus indices 0-23), columns B-F, I-M
# corresponding to slack bus voltage change from 1.05 pu to 1.02 pu (380 kV case)
$colLetters = @("B","C","D","E","F","I","J","K","L","M")

$data = @(
    @(1.02, 1.017506523434188, 1.02335253314526, 1.018887464434443, 1.015869691325769, 1.026927994660535, 1.022720250132159, 1.026184222918567, 1.021732345631859, 1.018723559631315),
    @(1.02, 1.018562034851111, 1.024119870061698, 1.019785563326101, 1.017568849712941, 1.027077371256606, 1.023411121460803, 1.026758528330028, 1.022436097497875, 1.020225499997718),
    @(1.02, 1.019244696683664, 1.024615847283845, 1.020366789292499, 1.018667787984145, 1.027172427672561, 1.023857311615084, 1.027128946393864, 1.022890942459054, 1.021196389025617),
    @(1.02, 1.019531611501204, 1.024824226797496, 1.020611160259838, 1.019129662561586, 1.027212006320759, 1.02404468744639, 1.027284384282853, 1.023082033227602, 1.021604326263554),
    @(1.02, 1.019579781335644, 1.024859207058205, 1.020652192620172, 1.019207206585259, 1.027218629287371, 1.024076136824634, 1.027310466237882, 1.023114110837523, 1.021672807756949),
    @(1.02, 1.019248530751122, 1.02461863216896, 1.020370054495975, 1.018673960032807, 1.027172958029569, 1.023859816134697, 1.027131024485271, 1.022893496317679, 1.021201840776527),
    @(1.02, 1.017863306583473, 1.023611970419874, 1.019190961664661, 1.016444044597401, 1.026978808537782, 1.022953909851761, 1.026378560019735, 1.021970291956028, 1.019231351786962),
    @(1.02, 1.01541981036163, 1.02183394875675, 1.017113956130374, 1.012510224125855, 1.026624434729114, 1.02135103634721, 1.025043441239661, 1.020339399011131, 1.015751382053041),
    @(1.02, 1.013788997663064, 1.020645788469741, 1.015729718600066, 1.009884199313228, 1.026379940321444, 1.020277987655714, 1.024147158011155, 1.019249338982973, 1.013425777247569),
    @(1.02, 1.013082384537353, 1.020130629594745, 1.015130420748529, 1.008746164484176, 1.026272114562841, 1.01981227181354, 1.023757578049253, 1.018776655610959, 1.012417332953745),
    @(1.02, 1.012819845721687, 1.019939174276089, 1.014907826787519, 1.008323295421758, 1.026231768917422, 1.019639120807948, 1.023612647056569, 1.018600976672673, 1.012042526721776),
    @(1.02, 1.012876164445788, 1.01998024673383, 1.014955573423391, 1.008414009317634, 1.026240436518577, 1.019676269702652, 1.023643745379771, 1.018638665093039, 1.012122934254481),
    @(1.02, 1.013060684462165, 1.020114805930144, 1.015112020824417, 1.008711213157949, 1.026268785585074, 1.019797962438672, 1.023745602575131, 1.01876213604855, 1.012386355995229),
    @(1.02, 1.013174363832497, 1.020197698709045, 1.015208414784361, 1.008894309785883, 1.026286213369844, 1.019872919656319, 1.023808330495616, 1.018838196859618, 1.012548628738959),
    @(1.02, 1.013835883372152, 1.020679963542554, 1.015769493754813, 1.009959706119626, 1.026387055092431, 1.020308872841237, 1.024172981829881, 1.01928069499752, 1.013492673269974),
    @(1.02, 1.014250712373064, 1.020982293713324, 1.016121466131068, 1.01062774024205, 1.026449786084819, 1.020582045017059, 1.02440132027778, 1.019558079512635, 1.014084455419502),
    @(1.02, 1.014492630823719, 1.021158572445819, 1.016326774043065, 1.011017302361961, 1.026486187099213, 1.020741277741338, 1.024534363158259, 1.019719807654928, 1.014429493371646),
    @(1.02, 1.014575111313388, 1.021218667871703, 1.016396780186255, 1.011150117855384, 1.026498566866866, 1.020795554388035, 1.024579703097855, 1.019774941685838, 1.014547119149053),
    @(1.02, 1.014206209743933, 1.020949863310586, 1.016083701967358, 1.010556075975648, 1.026443075184663, 1.020552746993917, 1.024376836528303, 1.019528325568332, 1.014020977190884),
    @(1.02, 1.013006349899533, 1.020075184431324, 1.0150659506425, 1.008623698315916, 1.026260445618759, 1.019762131482711, 1.023715614349976, 1.018725779789709, 1.012308791140859),
    @(1.02, 1.012251537196672, 1.019524646500894, 1.014426118595012, 1.007407850025688, 1.026143916006149, 1.019264094038661, 1.023298583362464, 1.018220589261873, 1.011230965108877),
    @(1.02, 1.01265171755005, 1.019816553440379, 1.014765299503109, 1.008052481429773, 1.026205852059711, 1.019528203243208, 1.02351978233495, 1.018488457386481, 1.011802467783703),
    @(1.02, 1.014226318703305, 1.020964517415652, 1.016100765941605, 1.010588458252884, 1.026446108135329, 1.020565985830769, 1.024387900122936, 1.019541770297612, 1.014049660697454),
    @(1.02, 1.016051825427843, 1.022294104261442, 1.017650832406737, 1.0135277869713, 1.026717501823305, 1.021766200210451, 1.025389692579326, 1.020761513477644, 1.016651996916855)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $colLetters.Count; $j++) {
        $cellRef = "$($colLetters[$j])$rowNum"
        $ws.Range($cellRef).Value = $rowVals[$j]
    }
}
